# Update field-sample rows for the new collection date (2020-09-22) across
# all the salt-route tabs, and leave the selection/active-tab state the way
# the workbook was left after this editing session (WIC active; DC + PBSF
# selections parked on their last data row).

$wb = $excel.ActiveWorkbook

# --- WIC: new row 6 ---
$ws = $wb.Worksheets.Item("WIC")
$ws.Range("A6").Value = 44096.388888888891
$ws.Range("B6").Value = 37.35

# --- YS: row 19 ---
$ws = $wb.Worksheets.Item("YS")
$ws.Range("A19").Value = 44096.490972222222
$ws.Range("B19").Value = 12.06
$ws.Range("C19").Value = 19.600000000000001

# --- SW: row 19 ---
$ws = $wb.Worksheets.Item("SW")
$ws.Range("A19").Value = 44096.511111111111
$ws.Range("B19").Value = 15.5
$ws.Range("C19").Value = 17.3

# --- YI: row 19 ---
$ws = $wb.Worksheets.Item("YI")
$ws.Range("A19").Value = 44096.402777777781
$ws.Range("B19").Value = 10.82
$ws.Range("C19").Value = 18.100000000000001

# --- YN: row 19 ---
$ws = $wb.Worksheets.Item("YN")
$ws.Range("A19").Value = 44096.418749999997
$ws.Range("B19").Value = 11.76
$ws.Range("C19").Value = 18.100000000000001

# --- 6MC: row 19 ---
$ws = $wb.Worksheets.Item("6MC")
$ws.Range("A19").Value = 44096.436111111114
$ws.Range("B19").Value = 10.039999999999999
$ws.Range("C19").Value = 15.4

# --- DC: row 19 (and leave selection parked on A19) ---
$ws = $wb.Worksheets.Item("DC")
$ws.Range("A19").Value = 44096.444444444445
$ws.Range("B19").Value = 8.8800000000000008
$ws.Range("C19").Value = 14.3
[void]$ws.Range("A19").Select()

# --- PBMS: row 19 ---
$ws = $wb.Worksheets.Item("PBMS")
$ws.Range("A19").Value = 44096.459722222222
$ws.Range("B19").Value = 39
$ws.Range("C19").Value = 16.7

# --- PBSF: row 19 (selection parked on A19, no longer the active tab) ---
$ws = $wb.Worksheets.Item("PBSF")
$ws.Range("A19").Value = 44096.468055555553
$ws.Range("B19").Value = 168.62
$ws.Range("C19").Value = 20.3
[void]$ws.Range("A19").Select()

# --- WIC becomes the active tab again, selection stays on A6 ---
$ws = $wb.Worksheets.Item("WIC")
[void]$ws.Activate()
[void]$ws.Range("A6").Select()
